# Update "想去人数" (want-to-go count) figures on the "展览" sheet and the
# corresponding rows on the aggregated "全部类型" sheet, per the latest
# gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 431   # 南宁·布谷鸟动漫展5th: 426 -> 431
$ws1.Range("F3").Value = 5295  # 南宁·2024良牙动漫秋季盛典（秋典）: 5264 -> 5295
$ws1.Range("F4").Value = 55    # 广西·THO04-永夜廻想: 52 -> 55
$ws1.Range("F5").Value = 60    # 南宁·花海演绎二次元水上派对: 57 -> 60
$ws1.Range("F7").Value = 508   # 南宁·万圣漫控嘉年华10: 507 -> 508

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 431   # 南宁·布谷鸟动漫展5th: 426 -> 431
$ws4.Range("F3").Value = 5295  # 南宁·2024良牙动漫秋季盛典（秋典）: 5264 -> 5295
$ws4.Range("F5").Value = 55    # 广西·THO04-永夜廻想: 52 -> 55
$ws4.Range("F6").Value = 60    # 南宁·花海演绎二次元水上派对: 57 -> 60
$ws4.Range("F9").Value = 508   # 南宁·万圣漫控嘉年华10: 507 -> 508
